$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 17 - shifts existing rows 17:52 down to 18:53,
# inheriting formatting (incl. the date number format in column D) from
# the surrounding rows.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly data point.
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(17, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(17, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(17, 4).Value = (Get-Date -Year 2021 -Month 8 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(17, 5).Value = 15
$ws.Cells.Item(17, 6).Value = 100112008
$ws.Cells.Item(17, 7).Value = "Coliflor"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Tercera"
$ws.Cells.Item(17, 10).Value = 1300
$ws.Cells.Item(17, 11).Value = 500
$ws.Cells.Item(17, 12).Value = 550
$ws.Cells.Item(17, 13).Value = 525
$ws.Cells.Item(17, 14).Value = "`$/unidad"
$ws.Cells.Item(17, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(17, 16).Value = 525
$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = "Hortaliza"
